$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update from the diff. Cells whose new value would be
# auto-recognized by Excel as a plain number (e.g. "557.21") are forced to
# remain text (matching the source workbook, where these are text/string
# values such as "557.21" stored as strings, not numeric cells), then the
# temporary number-format override is cleared by resetting the cell style
# back to "Normal" so no stray formatting is left behind.

$ws.Range("D2").Value = "68.318.88"
$ws.Range("D3").Value = "2.449.96"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("D9").Value = "2.451.09"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.08%  "
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "68.256.77"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "336.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +3.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "426.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.78%  "
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.72"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.481"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0917"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.49%  "
